$wb = $excel.ActiveWorkbook

# ---- Sheet 1 ----
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 1359
$ws.Range("G5").Value = "不可售"
$ws.Range("F6").Value = 40577
$ws.Range("F7").Value = 1442
$ws.Range("F8").Value = 267
$ws.Range("F11").Value = 5590
$ws.Range("F12").Value = 394
$ws.Range("F13").Value = 1033
$ws.Range("F14").Value = 2695
$ws.Range("F15").Value = 6300
$ws.Range("F16").Value = 166
$ws.Range("F17").Value = 1181
$ws.Range("F18").Value = 684
$ws.Range("F21").Value = 1091
$ws.Range("F23").Value = 61
$ws.Range("F25").Value = 135
$ws.Range("F27").Value = 854
$ws.Range("F31").Value = 1115
$ws.Range("F33").Value = 14
$ws.Range("F34").Value = 189
$ws.Range("F36").Value = 195
$ws.Range("F37").Value = 1143
$ws.Range("F39").Value = 65

# ---- Sheet 2 ----
$ws = $wb.Worksheets.Item(2)
$ws.Range("F5").Value = 503
$ws.Range("F10").Value = 26
$ws.Range("F16").Value = 540
$ws.Range("F27").Value = 106
$ws.Range("F28").Value = 531
$ws.Range("G28").Value = 399
$ws.Range("F29").Value = 929
$ws.Range("F30").Value = 542
$ws.Range("F32").Value = 71
$ws.Range("F35").Value = 89
$ws.Range("F36").Value = 115
$ws.Range("F38").Value = 42

# ---- Sheet 3 ----
$ws = $wb.Worksheets.Item(3)
$ws.Range("F5").Value = 813
$ws.Range("F6").Value = 507
$ws.Range("F7").Value = 274

# ---- Sheet 4 ----
$ws = $wb.Worksheets.Item(4)
$ws.Range("F4").Value = 1359
$ws.Range("B5").Value = "2024-04-14"
$ws.Range("C5").Value = "上海·青春无限乐团钢琴独奏&芭蕾舞表演梦幻联动表演经典影视作品曲目"
$ws.Range("D5").Value = "南苏州路1247号2楼 八号桥艺术空间"
$ws.Range("E5").Value = "2024.04.14 15:00-06.09 20:20"
$ws.Range("F5").Value = 11
$ws.Range("G5").Value = 238
$ws.Range("H5").Value = "https://show.bilibili.com/platform/detail.html?id=83604"
$ws.Range("I5").Value = "//i1.hdslb.com/bfs/openplatform/202403/WqNL1MyY1711595849962.jpeg"
$ws.Range("B6").Value = "2024-04-24"
$ws.Range("C6").Value = "上海·「NIJISANJI EN x animate cafe」"
$ws.Range("D6").Value = "西藏北路198号大悦城北座8楼N809-1 animate cafe上海店"
$ws.Range("E6").Value = "2024.04.24 00:00-05.22 23:59"
$ws.Range("F6").Value = 813
$ws.Range("G6").Value = 30
$ws.Range("H6").Value = "https://show.bilibili.com/platform/detail.html?id=83223"
$ws.Range("I6").Value = "//i1.hdslb.com/bfs/openplatform/202404/U4FVZbqr1713257119501.jpeg"
$ws.Range("C7").Value = "上海·国潮二次元沉浸式互动喜剧"
$ws.Range("D7").Value = "南京西路1038号梅龙镇广场（三楼307-308） 优+橘子喜剧空间"
$ws.Range("E7").Value = "2024.04.27 10:30-05.21 20:30"
$ws.Range("F7").Value = 15
$ws.Range("G7").Value = 78
$ws.Range("H7").Value = "https://show.bilibili.com/platform/detail.html?id=84771"
$ws.Range("I7").Value = "//i2.hdslb.com/bfs/openplatform/202404/Xw69nnAy1713779819613.jpeg"
$ws.Range("B8").Value = "2024-04-27"
$ws.Range("C8").Value = "上海·樱桃小丸子限定快闪-人累拯救包"
$ws.Range("D8").Value = "西藏北路166号 上海静安大悦城南座"
$ws.Range("E8").Value = "2024.04.27 10:00-05.26 22:00"
$ws.Range("F8").Value = 97
$ws.Range("G8").Value = 49
$ws.Range("H8").Value = "https://show.bilibili.com/platform/detail.html?id=84188"
$ws.Range("I8").Value = "//i1.hdslb.com/bfs/openplatform/202404/Kn9wkONJ1712893328433.png"
$ws.Range("C9").Value = "上海·2024《命运/冠位指定 冠位时间神殿所罗门》  萌果酱谷子咖啡"
$ws.Range("D9").Value = "南京东路340号百联ZX 萌果酱谷子咖啡（百联）"
$ws.Range("E9").Value = "2024.04.30 00:00-06.13 23:59"
$ws.Range("F9").Value = 507
$ws.Range("G9").Value = 30
$ws.Range("H9").Value = "https://show.bilibili.com/platform/detail.html?id=84585"
$ws.Range("I9").Value = "//i1.hdslb.com/bfs/openplatform/202404/QkbnfGA81713509547575.jpeg"
$ws.Range("F10").Value = 274
$ws.Range("B11").Value = "2024-04-30"
$ws.Range("C11").Value = "上海·T1 POP-UP STORE"
$ws.Range("D11").Value = "西藏北路166号（地铁8号线曲阜路下） 静安大悦城"
$ws.Range("E11").Value = "2024.04.30 00:00-05.16 23:59"
$ws.Range("F11").Value = 274
$ws.Range("G11").Value = 10
$ws.Range("H11").Value = "https://show.bilibili.com/platform/detail.html?id=84696"
$ws.Range("I11").Value = "//i2.hdslb.com/bfs/openplatform/202404/Azpxszvb1713853472117.png"
$ws.Range("C12").Value = "上海·街舞音乐剧《时光代理人：法则游戏》"
$ws.Range("D12").Value = "牛庄路704号 中国大戏院"
$ws.Range("E12").Value = "2024.05.01 19:30-05.19 21:00"
$ws.Range("F12").Value = 503
$ws.Range("G12").Value = 188
$ws.Range("H12").Value = "https://show.bilibili.com/platform/detail.html?id=82995"
$ws.Range("I12").Value = "//i1.hdslb.com/bfs/openplatform/202403/p9ZC2azX1710816437198.png"
$ws.Range("F13").Value = 1442
$ws.Range("F16").Value = 5590
$ws.Range("F17").Value = 394
$ws.Range("F18").Value = 1033
$ws.Range("F19").Value = 2695
$ws.Range("F20").Value = 26
$ws.Range("F21").Value = 6300
$ws.Range("F23").Value = 166
$ws.Range("F24").Value = 1181
$ws.Range("F26").Value = 540
$ws.Range("F27").Value = 684
$ws.Range("F29").Value = 1091
$ws.Range("F31").Value = 61
$ws.Range("F32").Value = 135
$ws.Range("F33").Value = 854
$ws.Range("F36").Value = 1115
$ws.Range("F39").Value = 929
$ws.Range("F40").Value = 542
$ws.Range("F41").Value = 189
$ws.Range("F43").Value = 71
$ws.Range("F44").Value = 195
$ws.Range("F46").Value = 89
$ws.Range("F47").Value = 115
$ws.Range("F49").Value = 65
$ws.Range("F50").Value = 42
